# Insert a new "slug" row right below the header row (row 2), pushing the
# existing measure/dim/type rows down by one (old row 2 -> 3, 3 -> 4, 4 -> 5),
# and fill the new row with the machine-friendly slug names for each column
# header (used to relate/relate-hierarchically the metadata columns, per
# issue #13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 2; existing rows 2.. shift down by one.
$ws.Rows.Item(2).Insert()

$slugs = @(
    "grandes-grupos-codigo",
    "municipio-codigo",
    "total",
    "sexo-codigo",
    "espanoles",
    "extranjeros",
    "sexo",
    "municipio-nombre",
    "grandes-grupos"
)

for ($i = 0; $i -lt $slugs.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $slugs[$i]
}
